$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename parameter labels in column A (time_* -> hr_*, latitude_* -> lat_*, longitude_* -> lon_*)
# Order matters for shared-string table ordering: hr_* first, then lat_*, then lon_*
$ws.Range("A11").Value = "hr_min"
$ws.Range("A12").Value = "hr_max"
$ws.Range("A13").Value = "hr_step"
$ws.Range("A2").Value = "lat_min"
$ws.Range("A3").Value = "lat_max"
$ws.Range("A4").Value = "lat_step"
$ws.Range("A5").Value = "lon_min"
$ws.Range("A6").Value = "lon_max"
$ws.Range("A7").Value = "lon_step"

# Update the view: remove frozen/scrolled topLeftCell (back to A1) and change selection to A9
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A9").Select()
